$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.009.37"
$ws.Range("E2").Value = "  +5.37%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.879.81"
$ws.Range("E3").Value = "  +4.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "281.42"
$ws.Range("E5").Value = "  +2.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5255"
$ws.Range("E7").Value = "  +3.92%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3531"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07049"
$ws.Range("E9").Value = "  +6.53%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.32"
$ws.Range("E10").Value = "  +1.73%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8173"
$ws.Range("E11").Value = "  -2.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07796"
$ws.Range("E12").Value = "  +0.49%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.881.70"
$ws.Range("E13").Value = "  +4.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.217"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "90.50"
$ws.Range("E15").Value = "  +3.41%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.62"
$ws.Range("E17").Value = "  +5.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008144"
$ws.Range("E18").Value = "  +2.26%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "27.038.73"
$ws.Range("E20").Value = "  +5.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.767"
$ws.Range("E21").Value = "  +0.98%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.19"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.243"
$ws.Range("E23").Value = "  +3.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.404"
$ws.Range("E24").Value = "  +14.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.79"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.61"
$ws.Range("E26").Value = "  +3.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.669"
$ws.Range("E27").Value = "  +0.98%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.76"
$ws.Range("E28").Value = "  +5.09%  "
$ws.Range("E29").Value = "  +1.83%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.387"
$ws.Range("E30").Value = "  +4.53%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08883"
$ws.Range("E31").Value = "  +1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04906"
$ws.Range("E32").Value = "  +2.18%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.175"
$ws.Range("E33").Value = "  +4.35%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7454"
$ws.Range("E34").Value = "  +3.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.895"
$ws.Range("E35").Value = "  +1.28%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.294"
$ws.Range("E36").Value = "  +8.94%  "
$ws.Range("E37").Value = "  +6.05%  "
$ws.Range("E38").Value = "  +2.69%  "
$ws.Range("E39").Value = "  +1.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9837"
$ws.Range("E40").Value = "  +3.04%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "117.13"
$ws.Range("E41").Value = "  +2.21%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.321"
$ws.Range("E42").Value = "  +2.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.186"
$ws.Range("E43").Value = "  +2.38%  "
$ws.Range("E44").Value = "  +0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4610"
$ws.Range("E45").Value = "  +0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1367"
$ws.Range("E46").Value = "  -0.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.502"
$ws.Range("E47").Value = "  +2.20%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "36.76"
$ws.Range("E48").Value = "  +2.50%  "
$ws.Range("E49").Value = "  +2.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05944"
$ws.Range("E50").Value = "  +2.56%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "61.84"
$ws.Range("E51").Value = "  +3.97%  "
